$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 169, shifting all
# following rows (old 169-260) down to (171-262).
$ws.Rows("169:170").Insert()

# Populate the first new row (169) with its data.
$ws.Range("A169").Value = 6
$ws.Range("B169").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C169").Value = "Metropolitana"
$ws.Range("D169").Value = 44873
$ws.Range("E169").Value = 13
$ws.Range("F169").Value = 100112001
$ws.Range("G169").Value = "Berenjena"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 1500
$ws.Range("K169").Value = 9000
$ws.Range("L169").Value = 10000
$ws.Range("M169").Value = 9467
$ws.Range("N169").Value = "`$/caja 50 unidades"
$ws.Range("O169").Value = "Región de Arica y Parinacota"
$ws.Range("P169").Value = 189
$ws.Range("Q169").Value = 50
$ws.Range("R169").Value = "Hortaliza"

# Populate the second new row (170) with its data.
$ws.Range("A170").Value = 6
$ws.Range("B170").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C170").Value = "Metropolitana"
$ws.Range("D170").Value = 44873
$ws.Range("E170").Value = 13
$ws.Range("F170").Value = 100112001
$ws.Range("G170").Value = "Berenjena"
$ws.Range("H170").Value = "Sin especificar"
$ws.Range("I170").Value = "Segunda"
$ws.Range("J170").Value = 600
$ws.Range("K170").Value = 7000
$ws.Range("L170").Value = 7000
$ws.Range("M170").Value = 7000
$ws.Range("N170").Value = "`$/caja 50 unidades"
$ws.Range("O170").Value = "Región de Arica y Parinacota"
$ws.Range("P170").Value = 140
$ws.Range("Q170").Value = 50
$ws.Range("R170").Value = "Hortaliza"
